$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9806.799999999999
$ws.Range("I43").Value = 5366.2
$ws.Range("J43").Value = 12027.1
$ws.Range("K43").Value = 5366.2
$ws.Range("L43").Value = 12027.1
$ws.Range("M43").Value = -5297.2
$ws.Range("N43").Value = -12165.1

$ws.Range("H76").Value = 3815.943
$ws.Range("I76").Value = 3334.373
$ws.Range("K76").Value = 3334.373
$ws.Range("M76").Value = -3019.373

$ws.Range("H79").Value = 3815.943
$ws.Range("I79").Value = 3334.373
$ws.Range("K79").Value = 3334.373
$ws.Range("M79").Value = -2242.373

$ws.Range("H116").Value = 6107.143
$ws.Range("I116").Value = 6192.3076
$ws.Range("K116").Value = 6192.3076
$ws.Range("M116").Value = -2750.3076

$ws.Range("H125").Value = 1596.6875
$ws.Range("J125").Value = 407
$ws.Range("L125").Value = 3663
$ws.Range("N125").Value = -8583

$ws.Range("H138").Value = 8625609
$ws.Range("I138").Value = 2159.8125
$ws.Range("J138").Value = 19239084
$ws.Range("K138").Value = 6479.4375
$ws.Range("L138").Value = 57717252
$ws.Range("M138").Value = -1339.4375
$ws.Range("N138").Value = -57727532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11130.932
$ws.Range("I32").Value = 12061.743
$ws.Range("J32").Value = 7511.1113
$ws.Range("K32").Value = 12061.743
$ws.Range("L32").Value = 7511.1113
$ws.Range("M32").Value = -11774.743
$ws.Range("N32").Value = -8085.1113

$ws.Range("H57").Value = 26985
$ws.Range("I57").Value = 26985
$ws.Range("K57").Value = 26985
$ws.Range("M57").Value = -26501

$ws.Range("H88").Value = 2732.6667
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2732.6667
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2732.6667
$ws.Range("N88").Value = -3544.6667
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 2732.6667
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2732.6667
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2732.6667
$ws.Range("N91").Value = -5540.6667
$ws.Range("M91").ClearContents()

$ws.Range("H122").Value = 4613.273
$ws.Range("I122").Value = 5461.76
$ws.Range("J122").Value = 1961.75
$ws.Range("K122").Value = 16385.28
$ws.Range("L122").Value = 5885.25
$ws.Range("M122").Value = -13935.28
$ws.Range("N122").Value = -10785.25

$ws.Range("H134").Value = 63269.555
$ws.Range("J134").Value = 63269.555
$ws.Range("L134").Value = 63269.555
$ws.Range("N134").Value = -73409.55499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17243214
$ws.Range("I86").Value = 1800.5
$ws.Range("J86").Value = 55557464
$ws.Range("K86").Value = 1800.5
$ws.Range("L86").Value = 55557464
$ws.Range("M86").Value = -677.5
$ws.Range("N86").Value = -55559710

$ws.Range("H89").Value = 17243214
$ws.Range("I89").Value = 1800.5
$ws.Range("J89").Value = 55557464
$ws.Range("K89").Value = 9002.5
$ws.Range("L89").Value = 277787320
$ws.Range("M89").Value = -3386.5
$ws.Range("N89").Value = -277798552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5750259
$ws.Range("I31").Value = 4719.9062
$ws.Range("J31").Value = 12821692
$ws.Range("K31").Value = 4719.9062
$ws.Range("L31").Value = 12821692
$ws.Range("M31").Value = -4424.9062
$ws.Range("N31").Value = -12822282

$ws.Range("H34").Value = 5750259
$ws.Range("I34").Value = 4719.9062
$ws.Range("J34").Value = 12821692
$ws.Range("K34").Value = 4719.9062
$ws.Range("L34").Value = 12821692
$ws.Range("M34").Value = -4517.9062
$ws.Range("N34").Value = -12822096

$ws.Range("H62").Value = 2875
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -3748

$ws.Range("H65").Value = 2875
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -18740

$ws.Range("H69").Value = 27910.3
$ws.Range("J69").Value = 34371.855
$ws.Range("L69").Value = 34371.855
$ws.Range("N69").Value = -35869.855

$ws.Range("H72").Value = 27910.3
$ws.Range("J72").Value = 34371.855
$ws.Range("L72").Value = 103115.565
$ws.Range("N72").Value = -110603.565

$ws.Range("H134").Value = 434603.3
$ws.Range("I134").Value = 1572.7174
$ws.Range("J134").Value = 2647871
$ws.Range("K134").Value = 4718.1522
$ws.Range("L134").Value = 7943613
$ws.Range("M134").Value = -2183.1522
$ws.Range("N134").Value = -7948683

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 50401.5
$ws.Range("I97").Value = 50401.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 151204.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -150708.5
$ws.Range("N97").ClearContents()

$ws.Range("H98").Value = 713
$ws.Range("I98").Value = 595
$ws.Range("J98").Value = 890
$ws.Range("K98").Value = 1785
$ws.Range("L98").Value = 2670
$ws.Range("M98").Value = -287
$ws.Range("N98").Value = -5666

$ws.Range("H107").Value = 860.0606
$ws.Range("J107").Value = 1766.1818
$ws.Range("L107").Value = 5298.5454
$ws.Range("N107").Value = -9138.545399999999

$ws.Range("H118").Value = 1928
$ws.Range("J118").Value = 1970.421
$ws.Range("L118").Value = 5911.263
$ws.Range("N118").Value = -8397.262999999999

$ws.Range("H131").Value = 1306.3219
$ws.Range("I131").Value = 2672.4443
$ws.Range("J131").Value = 949.942
$ws.Range("K131").Value = 8017.3329
$ws.Range("L131").Value = 2849.826
$ws.Range("M131").Value = -2977.3329
$ws.Range("N131").Value = -12929.826

$ws.Range("H132").Value = 4137.1816
$ws.Range("J132").Value = 5929.2856
$ws.Range("L132").Value = 53363.5704
$ws.Range("N132").Value = -58423.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 34164.8
$ws.Range("I70").Value = 152249.5
$ws.Range("K70").Value = 152249.5
$ws.Range("M70").Value = -151979.5

$ws.Range("H73").Value = 34164.8
$ws.Range("I73").Value = 152249.5
$ws.Range("K73").Value = 152249.5
$ws.Range("M73").Value = -151313.5

$ws.Range("H122").Value = 22225522
$ws.Range("I122").Value = 66666664
$ws.Range("J122").Value = 4950
$ws.Range("K122").Value = 199999992
$ws.Range("L122").Value = 14850
$ws.Range("M122").Value = -199997542
$ws.Range("N122").Value = -19750

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6506.8696
$ws.Range("I122").Value = 7618
$ws.Range("J122").Value = 5488.3335
$ws.Range("K122").Value = 22854
$ws.Range("L122").Value = 16465.0005
$ws.Range("M122").Value = -20404
$ws.Range("N122").Value = -21365.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8063.222
$ws.Range("J74").Value = 8500
$ws.Range("L74").Value = 8500
$ws.Range("N74").Value = -10372

$ws.Range("H77").Value = 8063.222
$ws.Range("J77").Value = 8500
$ws.Range("L77").Value = 25500
$ws.Range("N77").Value = -34860

$ws.Range("H113").Value = 1444.7778
$ws.Range("I113").Value = 1101.5385
$ws.Range("K113").Value = 3304.6155
$ws.Range("M113").Value = -1134.6155
